$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PowerConditioner")

# Update F10/F11/F12 values
$ws.Range("F10").Value = 53
$ws.Range("F11").Value = 43
$ws.Range("F12").Value = 33

# Update formulas for G10 and the shared formula G11:G12
$ws.Range("G10").Formula = "=`$B`$6*(1 + `$B`$5/F10)"
$ws.Range("G11").Formula = "=`$B`$6*(1 + `$B`$5/F11)"
$ws.Range("G12").Formula = "=`$B`$6*(1 + `$B`$5/F12)"

# Add new formatted (empty) cells K11 and K12 with the "0.0" number format (style index 1)
$ws.Range("K11:K12").NumberFormat = "0.0"

# Update the active selection to F12
$ws.Range("F12").Select()

$wb.Save()
